$wb = $excel.ActiveWorkbook

# ----- Sheet1 (Лист1): add new cells in row 1, 2, 3 -----
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("B1").Value = "9 749 ₽"
$ws1.Range("C1").Value = "977 ₽"
$ws1.Range("D1").Value = "6 913,27 ₽"
$ws1.Range("E1").Value = "750 ₽"
$ws1.Range("F1").Value = "10 181 ₽"
$ws1.Range("G1").Value = "1 729 ₽"
$ws1.Range("H1").Value = "11 835,01 ₽"
$ws1.Range("I1").Value = "8 092,97 ₽"
$ws1.Range("J1").Value = "6 710,31 ₽"
$ws1.Range("K1").Value = "6 152,18 ₽"

$ws1.Range("B2").Value = "329 ₽"
$ws1.Range("C2").Value = "449 ₽"
$ws1.Range("D2").Value = "146 ₽"

$ws1.Range("B3").Value = "Не найдено"

# ----- New sheet: AliexpressData (placed right after Лист1) -----
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "AliexpressData"

# Match Excel's default outline/page-margin settings for this sheet
$ws2.Outline.SummaryRow = 1
$ws2.Outline.SummaryColumn = 1
$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36

$ws2.Range("A1").Value = "АВР-Б-100-2-1"
$ws2.Range("B1").Value = "9 749 ₽"
$ws2.Range("C1").Value = "977 ₽"
$ws2.Range("D1").Value = "6 913,27 ₽"
$ws2.Range("E1").Value = "750 ₽"
$ws2.Range("F1").Value = "10 181 ₽"
$ws2.Range("G1").Value = "1 729 ₽"
$ws2.Range("H1").Value = "11 835,01 ₽"
$ws2.Range("I1").Value = "8 092,97 ₽"
$ws2.Range("J1").Value = "6 710,31 ₽"
$ws2.Range("K1").Value = "6 152,18 ₽"

$ws2.Range("A2").Value = "bababab"
$ws2.Range("B2").Value = "329 ₽"
$ws2.Range("C2").Value = "449 ₽"
$ws2.Range("D2").Value = "146 ₽"

$ws2.Range("A3").Value = "'267515"
$ws2.Range("B3").Value = "Не найдено"

# Restore Лист1 as the active/selected sheet
$ws1.Activate() | Out-Null
$ws1.Range("A3").Select() | Out-Null
